# Redesign Excel Overview and Adjust Code for Redesign
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (old rows 8 & 9: "Records Row" / "Records Banks Column")
# so the remaining rows shift up and the sheet ends at row 7.
$ws.Rows("8:9").Delete() | Out-Null

# Update labels (column A) to reflect the redesign
$ws.Range("A2").Value = "Current Income Row in IC Sheet"
$ws.Range("A3").Value = "Current Expense Row in IC Sheet"
$ws.Range("A6").Value = "Records Row"
$ws.Range("A7").Value = "Records Banks Column"

# Update values (column B) to reflect the redesign
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 6
$ws.Range("B7").Value = "J"

$wb.Save()
